$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '30.392.95'
$c = $ws.Range("E2")
$c.NumberFormat = "@"
$c.Value = '  +0.21%  '
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '1.938.37'
$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = '  +0.23%  '
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '1.000'
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '0.7705'
$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = '  +8.67%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '247.76'
$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = '  -1.17%  '
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.9996'
$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = '  -0.10%  '
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '28.08'
$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value = '  +1.48%  '
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.3218'
$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = '  -2.56%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '0.07117'
$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = '  -2.53%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.7857'
$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = '  -2.42%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.08032'
$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value = '  -0.64%  '
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '1.939.84'
$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = '  +0.28%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '5.390'
$c = $ws.Range("E14")
$c.NumberFormat = "@"
$c.Value = '  -1.55%  '
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '95.12'
$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = '  +0.56%  '
$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = '  -3.34%  '
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '30.394.43'
$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = '  +0.22%  '
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '256.22'
$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = '  +1.11%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '0.000008033'
$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = '  -1.85%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '5.843'
$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = '  +0.89%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '2.194.17'
$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = '  +0.24%  '
$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = '  -0.07%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '0.9996'
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '6.769'
$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = '  -2.98%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '9.633'
$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = '  -1.17%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '164.08'
$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = '  -0.76%  '
$c = $ws.Range("E27")
$c.NumberFormat = "@"
$c.Value = '  +4.86%  '
$c = $ws.Range("E28")
$c.NumberFormat = "@"
$c.Value = '  -0.89%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '2.301'
$c = $ws.Range("E29")
$c.NumberFormat = "@"
$c.Value = '  -1.96%  '
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '1.365'
$c = $ws.Range("E30")
$c.NumberFormat = "@"
$c.Value = '  +1.18%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '1.526'
$c = $ws.Range("E31")
$c.NumberFormat = "@"
$c.Value = '  -0.96%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '4.442'
$c = $ws.Range("E32")
$c.NumberFormat = "@"
$c.Value = '  +0.54%  '
$c = $ws.Range("E33")
$c.NumberFormat = "@"
$c.Value = '  -0.43%  '
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '0.05203'
$c = $ws.Range("E34")
$c.NumberFormat = "@"
$c.Value = '  +0.08%  '
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '1.287'
$c = $ws.Range("E35")
$c.NumberFormat = "@"
$c.Value = '  +1.91%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '0.7537'
$c = $ws.Range("E36")
$c.NumberFormat = "@"
$c.Value = '  +0.91%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '2.773'
$c = $ws.Range("E37")
$c.NumberFormat = "@"
$c.Value = '  -0.50%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.01979'
$c = $ws.Range("E38")
$c.NumberFormat = "@"
$c.Value = '  +0.66%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '2.812'
$c = $ws.Range("E39")
$c.NumberFormat = "@"
$c.Value = '  +0.12%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '79.12'
$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value = '  +0.16%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '6.498'
$c = $ws.Range("E41")
$c.NumberFormat = "@"
$c.Value = '  +1.28%  '
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '0.4532'
$c = $ws.Range("E42")
$c.NumberFormat = "@"
$c.Value = '  +0.18%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '1.986'
$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = '  -1.47%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '1.000'
$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = '  -0.01%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '0.8366'
$c = $ws.Range("E45")
$c.NumberFormat = "@"
$c.Value = '  -1.05%  '
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '101.53'
$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = '  -0.03%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '9.802'
$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value = '  +0.50%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '7.519'
$c = $ws.Range("E48")
$c.NumberFormat = "@"
$c.Value = '  +0.90%  '
$c = $ws.Range("E49")
$c.NumberFormat = "@"
$c.Value = '  +2.06%  '
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '985.19'
$c = $ws.Range("E50")
$c.NumberFormat = "@"
$c.Value = '  +11.41%  '
$c = $ws.Range("B51")
$c.NumberFormat = "@"
$c.Value = 'Algorand'
$c = $ws.Range("C51")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '0.1188'
$c = $ws.Range("E51")
$c.NumberFormat = "@"
$c.Value = '  +4.79%  '
